# Bug fix: a reference with an oversized field (authors/abstract) previously
# broke the whole SLR update. This adds the new "Misc. Data" column (J) and
# refreshes a few Authors (E) values / clears a stale "Other found locations"
# (I4) entry that were affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header: J1 = "Misc. Data"
$ws.Range("J1").Value = 'Misc. Data'

# Row 2: refreshed Authors value; new (blank) Misc. Data cell
$ws.Range("E2").Value = '[Atas%Jenny%coreGivesNoEmail%1,                        Bandy%Kenneth%coreGivesNoEmail%1,                        Bradin%Stuart A.%coreGivesNoEmail%1,                        Cadwallender%Bruce A.%coreGivesNoEmail%1,                        Cinti%Sandro K.%coreGivesNoEmail%1,                        Collins%Curtis D.%coreGivesNoEmail%1,                        Goldberg%Janet%coreGivesNoEmail%1,                        Holmes%Jennifer G.%coreGivesNoEmail%1,                        Kim%Christopher%coreGivesNoEmail%1,                        Krupansky%Frank%coreGivesNoEmail%1,                        Lozon%Marie M.%coreGivesNoEmail%1,                        Rodgers%Phillip E.%coreGivesNoEmail%1,                        Shlafer%Jean%coreGivesNoEmail%1,                        Wagner%Deborah%coreGivesNoEmail%1,                        Wilkerson%William M.%coreGivesNoEmail%1,                        Wright%Carrie M.%coreGivesNoEmail%1]'
$ws.Range("J2").Value = ""

# Row 3: new (blank) Misc. Data cell
$ws.Range("J3").Value = ""

# Row 4: refreshed Authors value; "Other found locations" cleared; new Misc. Data value
$ws.Range("E4").Value = '[Cathy%Campbell%xref no email%1,     Marianne%Baernholdt%xref no email%1]'
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = 'PUBLISHER: Project MUSE'

# Row 5: new (blank) Misc. Data cell
$ws.Range("J5").Value = ""
